$wb = $excel.ActiveWorkbook

# 1. Reorder worksheets: move "review_info" to be the first sheet (before "hotel_info").
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))

# 2. Insert a new "State" column into "hotel_info", between "Hotel_Name" and "City",
#    and populate it with "Louisiana" for the existing data row.
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Columns.Item(3).EntireColumn.Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"
